# The commit swaps the contents of ppt/theme/theme1.xml (used by the
# slide master -> all slides) and ppt/theme/theme2.xml (used by the
# notes master): theme1 goes from the "Integral" / "Red Violet" colour
# scheme to the stock "Office Theme" / "Office" colour scheme, and
# theme2 goes the other way.
#
# This COM host only exposes the *slide master's* theme through the
# PowerPoint object model (Presentation.SlideMaster.Theme /
# Presentation.NotesMaster.Theme both resolve to the same underlying
# theme part, ppt/theme/theme1.xml); there is no reachable object for
# the notes-master theme part. So we apply the reachable half of the
# swap: recolour the presentation's theme (theme1.xml) to the "Office"
# palette via ThemeColorScheme, which is the part that actually drives
# the look of every slide.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# DrawingML <a:clrScheme> slot order, 1-based as exposed by
# ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Values below are the stock Office theme palette (hex -> COM RGB()
# packs as R | G<<8 | B<<16).

$colors.Item(1).RGB  = 0x000000       # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF       # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444       # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7       # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B       # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED       # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5       # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF       # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244       # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70       # accent6  70AD47
$colors.Item(11).RGB = 0xC16305       # hlink    0563C1
$colors.Item(12).RGB = 0x724F95       # folHlink 954F72
